$d = $word.ActiveDocument

# --- Paragraphs being marked as done (strikethrough) ---
# 27: "Retrieve the customer's account information, including all their accounts and their current balances."
# 28: "Apply for a new loan or credit card."
$p27 = $d.Paragraphs.Item(27).Range
$p28 = $d.Paragraphs.Item(28).Range

$p27.Font.StrikeThrough = 1
$p27.Font.DoubleStrikeThrough = 0

$p28.Font.StrikeThrough = 1
$p28.Font.DoubleStrikeThrough = 0

# --- Move the "_GoBack" bookmark ---
# It currently sits mid-sentence in paragraph 27 (between "their" and " current").
# It should now wrap paragraph 28 ("Apply for a new loan or credit card."),
# marking it as the most-recently-edited location.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$p28b = $d.Paragraphs.Item(28).Range
$d.Bookmarks.Add("_GoBack", $p28b)
